# Refresh crypto price/volume snapshot (Price column D, Volume(1h) column E)
# for the coin rows whose figures moved in this data pull.
# Values are written as literal text (NumberFormat "@") so the figures
# round-trip exactly as they appeared in the source feed, matching the
# existing inline-string cells for these columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.07%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.41%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.133"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.14%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07606"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.92%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.623"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.17%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.455"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.98%"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9014"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.21%"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1126"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "12.20%"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1766"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.06%"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09204"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.16%"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04184"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.48%"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1048"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.69%"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001253"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.45%"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005873"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.92%"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.357"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.06%"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.235"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.05%"

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.95%"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.552"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-7.15%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1360"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.37%"

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-15.08%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04131"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.13%"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001226"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.56%"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004000"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.56%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001301"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.50%"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02400"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "2.07%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05184"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.04%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007764"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.51%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1300"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.57%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006956"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.17%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001971"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.81%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007597"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.13%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3054"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.08%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006738"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.52%"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.06%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03142"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "823.25%"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004202"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-15.99%"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.06%"

